# Auto-generated Excel COM-interop edit script
# Updates cryptos list values (Price / Volume(1h) / swapped rows) per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.488.19"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.661.85"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.12"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.63"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.661.09"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.99"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.151.33"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.378.67"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000185"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.25"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.656.51"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.39"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.17"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.05"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.01"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.33"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0971"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.13"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "496.55"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.94"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.52"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.94"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.38"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.76"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.01"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.67"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.28"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.559"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  -2.88%  "
